# "Generate Report for Handback" - update the localization-status report to
# reflect that the handback has completed and is now in sync with en-US.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: every cell that showed "Ready for handoff" (the Overview
# summary columns plus each language sheet's own Status column) now reports
# that the handback has completed and is in sync with en-US.
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"

$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

# Widen the status columns now that the text is longer.
$wsOverview.Columns("E").ColumnWidth = 29.166666666666668
$wsOverview.Columns("F").ColumnWidth = 29.166666666666668

# --- zh-cn detail sheet: refresh the handback timestamp and clear the
# stale "handback not latest" error now that the handback is in sync.
$wsZhCn.Range("K2").Value = "2016-08-13 06:44:41"
$wsZhCn.Range("K3").Value = "2016-08-13 06:44:41"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Columns("C").ColumnWidth = 29.166666666666668
$wsZhCn.Columns("P").ColumnWidth = 12.833333333333334

# --- de-de detail sheet: same treatment - refresh handback timestamp and
# clear the error now that the handback is in sync.
$wsDeDe.Range("K2").Value = "2016-08-13 06:44:51"
$wsDeDe.Range("K3").Value = "2016-08-13 06:44:51"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Columns("C").ColumnWidth = 29.166666666666668
$wsDeDe.Columns("P").ColumnWidth = 12.833333333333334
